$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vm_pu values for the 380 kV case (rows 2-25), columns B,C,D,E,F,I,J,K,L,M
$columns = @("B", "C", "D", "E", "F", "I", "J", "K", "L", "M")

$newValues = @{
    2 = @(1.02, 1.038927444142529, 1.047616587205179, 1.03752478985291, 1.056188527193665, 1.042043491777227, 1.044022227579138, 1.050378726037938, 1.040315426889857, 1.058926923289173)
    3 = @(1.02, 1.039787584377955, 1.048301733414593, 1.03825420163419, 1.057004621183756, 1.042256228456884, 1.04452781047057, 1.05087608095883, 1.040854877054477, 1.059556600591692)
    4 = @(1.02, 1.040344838241483, 1.048745672616952, 1.038727138843579, 1.057533635725298, 1.042393072671594, 1.044854980474361, 1.051197826453084, 1.041204226393943, 1.059964336641698)
    5 = @(1.02, 1.040579270305734, 1.048932447431524, 1.03892618944282, 1.057756258273463, 1.042450407178453, 1.044992527229103, 1.05133306863562, 1.041351160674433, 1.06013581737374)
    6 = @(1.02, 1.040618641965935, 1.048963816065456, 1.038959624220237, 1.057793650707587, 1.042460022449987, 1.045015622160384, 1.051355775203685, 1.041375835523683, 1.060164613713851)
    7 = @(1.02, 1.040347970097522, 1.048748167752391, 1.038729797674598, 1.057536609537641, 1.042393839544426, 1.044856818364805, 1.051199633644523, 1.041206189472232, 1.059966627708473)
    8 = @(1.02, 1.039217989639331, 1.04784800948885, 1.037771098331759, 1.056464132710305, 1.042115554542539, 1.044193085933898, 1.050546824478412, 1.040497675932354, 1.059139664117746)
    9 = @(1.02, 1.037232141463659, 1.046266518410762, 1.036089171781896, 1.054581630431146, 1.041619004688498, 1.043023745370418, 1.049395965223944, 1.039251458697216, 1.057684760637538)
    10 = @(1.02, 1.035911915938126, 1.045215458859367, 1.034972986069888, 1.053331677643743, 1.041283862718402, 1.042244423016686, 1.048628451605641, 1.038422258061563, 1.056716473988286)
    11 = @(1.02, 1.035341134963285, 1.044761135949633, 1.034490897304001, 1.052791655963897, 1.041137777226265, 1.041907040342505, 1.048296061516011, 1.038063604384843, 1.05629760625874)
    12 = @(1.02, 1.035129256117222, 1.044592501117123, 1.034312014254795, 1.052591252536392, 1.04108337004742, 1.041781733058564, 1.048172590630502, 1.037930445279824, 1.056142082857525)
    13 = @(1.02, 1.035174698731244, 1.044628668359924, 1.034350376821384, 1.052634231363757, 1.041095047100193, 1.041808611347836, 1.048199075822803, 1.037959005584362, 1.056175440290797)
    14 = @(1.02, 1.035323618236243, 1.044747194064455, 1.034476106966637, 1.052775086775541, 1.041133282855881, 1.041896682157937, 1.04828585549698, 1.03805259616251, 1.056284749357027)
    15 = @(1.02, 1.035415390344962, 1.044820237740815, 1.03455359818483, 1.05286189697259, 1.041156822038105, 1.041950947074821, 1.048339322509792, 1.03811026851013, 1.056352106640029)
    16 = @(1.02, 1.035949815582414, 1.045245627633101, 1.035005006722288, 1.053367542893764, 1.04129353762407, 1.042266815558441, 1.048650510268662, 1.03844606918025, 1.056744281567745)
    17 = @(1.02, 1.036285284588227, 1.045512677012955, 1.035288493176371, 1.053685048195894, 1.041379037371464, 1.042464970909186, 1.048845697258795, 1.038656814985568, 1.056990392581254)
    18 = @(1.02, 1.036481043303233, 1.045668518771379, 1.03545396429457, 1.053870360873912, 1.041428814595642, 1.042580558148325, 1.0489595414559, 1.038779777531209, 1.057133984111889)
    19 = @(1.02, 1.036547806389196, 1.045721669682722, 1.035510405673685, 1.053933567528936, 1.041445771489859, 1.04261997148594, 1.048998358483767, 1.038821710998733, 1.057182951703015)
    20 = @(1.02, 1.036249283102182, 1.045484017229095, 1.035258065516693, 1.053650970744404, 1.041369873706232, 1.042443710025929, 1.04882475605686, 1.038634200005979, 1.056963983138597)
    21 = @(1.02, 1.035279761424626, 1.044712287845966, 1.034439077418268, 1.052733603253996, 1.041122027362495, 1.041870747167146, 1.048260301211869, 1.038025034379259, 1.05625255881707)
    22 = @(1.02, 1.034670964343533, 1.044227771484142, 1.033925225204375, 1.052157887309715, 1.040965360728593, 1.041510570795771, 1.047905369650113, 1.037642380461291, 1.055805621806513)
    23 = @(1.02, 1.034993624754437, 1.044484555767248, 1.034197525179495, 1.052462983202845, 1.041048491655581, 1.04170150021728, 1.048093528585626, 1.03784519864674, 1.056042516560257)
    24 = @(1.02, 1.036265550360625, 1.045496967116029, 1.035271814101175, 1.053666368515987, 1.041374014660574, 1.04245331688375, 1.04883421850045, 1.03864441862431, 1.056975916308219)
    25 = @(1.02, 1.037744890366884, 1.046674803537988, 1.036523099271999, 1.055067420901252, 1.041748102131622, 1.04332601155977, 1.049693544073044, 1.039573357792826, 1.058060604336421)
}

foreach ($row in $newValues.Keys) {
    $rowValues = $newValues[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $cellRef = "{0}{1}" -f $columns[$i], $row
        $ws.Range($cellRef).Value = $rowValues[$i]
    }
}
